# Preparation of publication 0.2.0
# - Bump the Version and Date metadata values
# - Insert a new "Jurisdiction" row on the Metadata sheet

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Update Version (row 3) and Date (row 8) values
$ws.Range("B3").Value = "0.2.0"
$ws.Range("B8").Value = "2023-10-19T17:05:12+00:00"

# Insert a new row right after "Contact" (row 10) for the new "Jurisdiction" property,
# copying the formatting of the existing data rows so the new cells keep the same style.
$ws.Rows.Item(11).Insert()
$ws.Range("A10:B10").Copy()
$ws.Range("A11:B11").PasteSpecial(-4122)

$ws.Range("A11").Value = "Jurisdiction"
$ws.Range("B11").Value = "iso:code:3166:FR"
